$d = $word.ActiveDocument

# 1. Fix the typo: "needs" -> "have" and remove "not" before "fixed."
$d.Content.Find.Execute(
    "needs to be set manually by the administrator until the problem is not fixed.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "have to be set manually by the administrator until the problem is fixed.", 2
) | Out-Null

# 2. Force a run boundary around "have" (mirrors how Word splits runs at edit points)
$rngHave = $d.Content
$rngHave.Find.Execute("have", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("tmpHaveStart", $d.Range($rngHave.Start, $rngHave.Start)) | Out-Null
$d.Bookmarks.Add("tmpHaveEnd", $d.Range($rngHave.End, $rngHave.End)) | Out-Null
$d.Bookmarks("tmpHaveStart").Delete()
$d.Bookmarks("tmpHaveEnd").Delete()

# 3. Force run boundaries around "django" and "oauth" (mirrors Word's spell-check run splitting)
$rngDjango = $d.Content
$rngDjango.Find.Execute("django", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("tmpDjangoStart", $d.Range($rngDjango.Start, $rngDjango.Start)) | Out-Null
$d.Bookmarks.Add("tmpDjangoEnd", $d.Range($rngDjango.End, $rngDjango.End)) | Out-Null
$d.Bookmarks("tmpDjangoStart").Delete()
$d.Bookmarks("tmpDjangoEnd").Delete()

$rngOauth = $d.Content
$rngOauth.Find.Execute("oauth", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("tmpOauthStart", $d.Range($rngOauth.Start, $rngOauth.Start)) | Out-Null
$d.Bookmarks.Add("tmpOauthEnd", $d.Range($rngOauth.End, $rngOauth.End)) | Out-Null
$d.Bookmarks("tmpOauthStart").Delete()
$d.Bookmarks("tmpOauthEnd").Delete()

# 4. Place the _GoBack bookmark right before "fixed." (the last edit position),
#    which also matches Word's automatic behavior of relocating this bookmark.
$rngFixed = $d.Content
$rngFixed.Find.Execute("fixed.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $d.Range($rngFixed.Start, $rngFixed.Start)) | Out-Null

Write-Output "Done"
